# Update countries & provincias Spain
# - Refresh COVID counters for a number of countries (rows identified by their
#   fixed worksheet row, data laid out País/Casos totales/Nuevos casos/
#   Casos activos/Recuperados/Casos críticos/Muertes hoy/Muertes in A:H).
# - "Georgia" moves above "Trinidad yTobago" in the source list, and
#   "Montserrat" moves above "Islas Malvinas" - since both rows keep their
#   place in the sheet, this shows up as the two rows swapping which country
#   (and data) they carry.
# - Refresh the "Datos actualizados ..." timestamp banner in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain counter refreshes (country identity / row unchanged) ----------

# Estados Unidos (row 4)
$ws.Range("B4").Value = 7098291
$ws.Range("C4").Value = 354
$ws.Range("D4").Value = 4347172
$ws.Range("E4").Value = 2545641
$ws.Range("G4").Value = 7
$ws.Range("H4").Value = 205478

# India (row 5)
$ws.Range("B5").Value = 5650540
$ws.Range("C5").Value = 10044
$ws.Range("E5").Value = 972850
$ws.Range("G5").Value = 56
$ws.Range("H5").Value = 90077

# Barein (row 56)
$ws.Range("E56").Value = 6806
$ws.Range("G56").Value = 2
$ws.Range("H56").Value = 229

# Singapur (row 57)
$ws.Range("B57").Value = 57639
$ws.Range("C57").Value = 12
$ws.Range("E57").Value = 350

# Armenia (row 63)
$ws.Range("B63").Value = 47877
$ws.Range("C63").Value = 210
$ws.Range("D63").Value = 43026
$ws.Range("E63").Value = 3909
$ws.Range("G63").Value = 4
$ws.Range("H63").Value = 942

# Australia (row 78)
$ws.Range("B78").Value = 26974
$ws.Range("C78").Value = 32
$ws.Range("D78").Value = 24416
$ws.Range("E78").Value = 1699

# Hungria (row 83)
$ws.Range("B83").Value = 20450
$ws.Range("C83").Value = 951
$ws.Range("D83").Value = 4644
$ws.Range("E83").Value = 15104
$ws.Range("G83").Value = 8
$ws.Range("H83").Value = 702

# Togo (row 161)
$ws.Range("B161").Value = 1572
$ws.Range("C161").Value = 12
$ws.Range("E161").Value = 288

# Taiwan (row 176)
$ws.Range("D176").Value = 480
$ws.Range("E176").Value = 22

# --- Reordered pair: Georgia now listed before Trinidad yTobago ----------
# Row 131 keeps showing fresh "Georgia" data, row 132 now carries the data
# that used to belong to "Trinidad yTobago" in row 131.
$ws.Range("A131").Value = "Georgia"
$ws.Range("B131").Value = 4140
$ws.Range("C131").Value = 227
$ws.Range("D131").Value = 1643
$ws.Range("E131").Value = 2473
$ws.Range("F131").Value = 0
$ws.Range("G131").Value = 1
$ws.Range("H131").Value = 24

$ws.Range("A132").Value = "Trinidad yTobago"
$ws.Range("B132").Value = 4026
$ws.Range("C132").Value = 0
$ws.Range("D132").Value = 1871
$ws.Range("E132").Value = 2090
$ws.Range("F132").Value = 0
$ws.Range("G132").Value = 0
$ws.Range("H132").Value = 65

# --- Reordered pair: Montserrat now listed before Islas Malvinas ---------
# Rows 214/215 swap their country + data wholesale (no new counters, a pure
# reorder).
$ws.Range("A214").Value = "Montserrat"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 12
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 1

$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("B215").Value = 13
$ws.Range("C215").Value = 0
$ws.Range("D215").Value = 13
$ws.Range("E215").Value = 0
$ws.Range("F215").Value = 0
$ws.Range("G215").Value = 0
$ws.Range("H215").Value = 0

# --- Timestamp banner ------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 23 de Septiembre de 2020 a las 09:34"
